# Logged Week 15 and simulated Week 16
# Appends the new week's per-play yardage samples to the YDS and ST
# "raw list" cells, and updates the season-to-date totals on the
# OFF / DEF / ST / TURNS / PEN summary sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: space-separated per-play yardage lists (Rush / Pass, Off / Def)
# ---------------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value() + " 1 6 5 4 5 5 28 4 5 0 5 1 -3 0 6 6 3 5 1 3 2 0"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value() + " 3 -1 7 8 4 13 17 2 3 2 0 2 3 10 7 9 12 2 5 0 0 2 1 -2 1 12 3 1"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value() + " 8 2 16 4 15 0 3 4 9 15 9 7 7 8 14 7 7 10 36"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value() + " 9 11 4 2 19 4 9 9 6 5 11 3 3 6 7 15 9 3 7 3 10 16 10 1 11 7 9 7"

# ---------------------------------------------------------------------------
# OFF sheet: row 2 = RATT (Home), row 3 = PATT (Road)
# ---------------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 145
$offWs.Range("F2").Value = 67
$offWs.Range("H2").Value = 7
$offWs.Range("J2").Value = 22
$offWs.Range("L2").Value = 242
$offWs.Range("M2").Value = 152
$offWs.Range("O2").Value = 23
$offWs.Range("P2").Value = 10
$offWs.Range("Q2").Value = 435

$offWs.Range("B3").Value = 7
$offWs.Range("C3").Value = 184
$offWs.Range("E3").Value = 25
$offWs.Range("F3").Value = 95
$offWs.Range("G3").Value = 25
$offWs.Range("H3").Value = 23
$offWs.Range("I3").Value = 66
$offWs.Range("J3").Value = 48
$offWs.Range("N3").Value = 15

# ---------------------------------------------------------------------------
# DEF sheet: row 2 = RATT (Home), row 3 = PATT (Road)
# ---------------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value = 186
$defWs.Range("D2").Value = 11
$defWs.Range("E2").Value = 5
$defWs.Range("F2").Value = 55
$defWs.Range("G2").Value = 66
$defWs.Range("I2").Value = 6
$defWs.Range("J2").Value = 25
$defWs.Range("L2").Value = 255
$defWs.Range("M2").Value = 169
$defWs.Range("Q2").Value = 452

$defWs.Range("C3").Value = 169
$defWs.Range("E3").Value = 23
$defWs.Range("F3").Value = 101
$defWs.Range("G3").Value = 35
$defWs.Range("H3").Value = 14
$defWs.Range("I3").Value = 58
$defWs.Range("J3").Value = 57
$defWs.Range("N3").Value = 18

# ---------------------------------------------------------------------------
# ST sheet: row 2 = "#" totals, row 3 = "TB" totals, rows 3-6 also hold the
# per-kick distance raw lists (D / RA / RM, paired columns B & D)
# ---------------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 58
$stWs.Range("D2").Value = 56
$stWs.Range("J2").Value = 176
$stWs.Range("K2").Value = 159
$stWs.Range("L2").Value = 110
$stWs.Range("M2").Value = 92

$stWs.Range("B3").Value = 26
$stWs.Range("D3").Value = $stWs.Range("D3").Value() + " 50 51"

$stWs.Range("B4").Value = $stWs.Range("B4").Value() + " 62 60"
$stWs.Range("D4").Value = $stWs.Range("D4").Value() + " 0 0"

$stWs.Range("B5").Value = $stWs.Range("B5").Value() + " 21 29"
$stWs.Range("D5").Value = $stWs.Range("D5").Value() + " 0 0 0 0"

$stWs.Range("B6").Value = $stWs.Range("B6").Value() + " 15"

# ---------------------------------------------------------------------------
# TURNS sheet: row 2 = Home, row 3 = Road
# ---------------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B2").Value = 6
$turnsWs.Range("D2").Value = 7
$turnsWs.Range("E2").Value = 7

$turnsWs.Range("D3").Value = 6

# ---------------------------------------------------------------------------
# PEN sheet: row 2 = False start, row 3 = Holding
# ---------------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value = 16
$penWs.Range("B3").Value = 10
